$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are written as text (avoid Excel auto-numeric coercion)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.252.97"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.562.74"
$ws.Range("E3").Value = "  +4.68%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "606.46"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "145.22"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("D7").Value = "3.561.91"
$ws.Range("E7").Value = "  +4.76%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").Value = "0.137"
$ws.Range("E10").Value = "  +1.62%  "
$ws.Range("D11").Value = "8.02"
$ws.Range("E11").Value = "  +1.32%  "
$ws.Range("E12").Value = "  +1.78%  "
$ws.Range("D13").Value = "4.166.35"
$ws.Range("E13").Value = "  +4.63%  "
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "30.15"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").Value = "3.563.28"
$ws.Range("E16").Value = "  +5.34%  "
$ws.Range("D17").Value = "66.341.03"
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  +10.40%  "
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("E21").Value = "  +1.57%  "
$ws.Range("D22").Value = "431.90"
$ws.Range("E22").Value = "  +3.75%  "
$ws.Range("D23").Value = "0.609"
$ws.Range("E23").Value = "  +5.22%  "
$ws.Range("D24").Value = "78.77"
$ws.Range("E24").Value = "  +1.87%  "
$ws.Range("D25").Value = "3.704.30"
$ws.Range("E25").Value = "  +4.52%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  +7.91%  "
$ws.Range("D28").Value = "8.06"
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("D29").Value = "2.52"
$ws.Range("E29").Value = "  +4.03%  "
$ws.Range("D30").Value = "9.19"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "3.555.56"
$ws.Range("E34").Value = "  +4.56%  "
$ws.Range("D35").Value = "25.43"
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").Value = "7.90"
$ws.Range("E38").Value = "  +4.61%  "
$ws.Range("D39").Value = "5.66"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D41").Value = "171.53"
$ws.Range("E41").Value = "  +1.01%  "
$ws.Range("D42").Value = "0.0855"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "5.22"
$ws.Range("E43").Value = "  +3.16%  "
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").Value = "46.07"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("E47").Value = "  +4.19%  "
$ws.Range("D48").Value = "26.03"
$ws.Range("E48").Value = "  -3.14%  "
$ws.Range("E49").Value = "  +4.06%  "
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("E51").Value = "  +3.71%  "

# Restore default style so no stray formatting is left on the cells
$ws.Range("D2:D51").Style = "Normal"

